$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly reported totals for "Spint( 40) - Day 7 - Test Case Summary"
$ws.Range("C39").Value = 2829
$ws.Range("C40").Value = 1367
$ws.Range("C41").Value = 803

# Reflect where the user was working when they saved (view/selection)
$ws.Range("C41").Select()
